$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2: update existing product row with new sku/id and new product name,
# switch status/visibility columns from text to numeric codes, and turn the
# image/small_image/thumbnail cells into hyperlinks (same URL text as before).
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = 11
$ws.Range("B2").Value = 6972
$ws.Range("C2").Value = "ELEGANCE GOLD MATIC EYE CONTOUR BLACK"
$ws.Range("D2").Value = "ELEGANCE GOLD MATIC EYE CONTOUR BLACK"
$ws.Range("I2").Value = "ELEGANCE GOLD MATIC EYE CONTOUR BLACK"
$ws.Range("J2").Value = "ELEGANCE GOLD MATIC EYE CONTOUR BLACK"
$ws.Range("U2").Value = 1
$ws.Range("V2").Value = 4

# ---------------------------------------------------------------------------
# Row 3: brand new product row
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 6974
$ws.Range("C3").Value = "ELEGANCE GOLD MATIC EYE CONTOUR BROWN"
$ws.Range("D3").Value = "ELEGANCE GOLD MATIC EYE CONTOUR BROWN"
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = "ELEGANCE"
$ws.Range("G3").Value = "MODERN PHARMACEUTICAL CO."
$ws.Range("H3").Value = "Normal"
$ws.Range("I3").Value = "ELEGANCE GOLD MATIC EYE CONTOUR BROWN"
$ws.Range("J3").Value = "ELEGANCE GOLD MATIC EYE CONTOUR BROWN"
$ws.Range("K3").Value = "media/import/pdc.png"
$ws.Range("L3").Value = "media/import/pdc.png"
$ws.Range("M3").Value = "media/import/pdc.png"
$ws.Range("N3").Value = 100000
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = "admin"
$ws.Range("S3").Value = "pdc"
$ws.Range("T3").Value = 1
$ws.Range("U3").Value = 1
$ws.Range("V3").Value = 4
$ws.Range("W3").Value = "Shipping"
$ws.Range("X3").Value = "simple"
$ws.Range("Y3").Value = "Default"

# ---------------------------------------------------------------------------
# Row 4: brand new product row
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = 19
$ws.Range("B4").Value = 7156
$ws.Range("C4").Value = "ELEGANCE GOLD SOFT SHINY EYE LINER BLACK"
$ws.Range("D4").Value = "ELEGANCE GOLD SOFT SHINY EYE LINER BLACK"
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = "ELEGANCE"
$ws.Range("G4").Value = "MODERN PHARMACEUTICAL CO."
$ws.Range("H4").Value = "Normal"
$ws.Range("I4").Value = "ELEGANCE GOLD SOFT SHINY EYE LINER BLACK"
$ws.Range("J4").Value = "ELEGANCE GOLD SOFT SHINY EYE LINER BLACK"
$ws.Range("K4").Value = "/home/pdcorders.com/media/import/pdc.png"
$ws.Range("L4").Value = "/home/pdcorders.com/media/import/pdc.png"
$ws.Range("M4").Value = "/home/pdcorders.com/media/import/pdc.png"
$ws.Range("N4").Value = 100000
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = "admin"
$ws.Range("S4").Value = "pdc"
$ws.Range("T4").Value = 1
$ws.Range("U4").Value = 1
$ws.Range("V4").Value = 4
$ws.Range("W4").Value = "Shipping"
$ws.Range("X4").Value = "simple"
$ws.Range("Y4").Value = "Default"

# ---------------------------------------------------------------------------
# Turn K2/L2/M2 into hyperlinks pointing at the same pdc.png URL they already
# display (this also introduces the "Hyperlink" cell style / font).
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("K2"), "https://www.pdcorders.com/media/import/pdc.png") | Out-Null
$ws.Hyperlinks.Add($ws.Range("L2"), "https://www.pdcorders.com/media/import/pdc.png") | Out-Null
$ws.Hyperlinks.Add($ws.Range("M2"), "https://www.pdcorders.com/media/import/pdc.png") | Out-Null

# ---------------------------------------------------------------------------
# Scroll the view over and select V4, matching the final sheetView state.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 8
$win.ScrollRow = 1
$ws.Range("V4").Select() | Out-Null
